$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the first 11 data rows (rows 2-12): Serial -> 526-622-584, pdv -> 1360
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = "526-622-584"
    $ws.Cells.Item($r, 2).Value = 1360
    $ws.Cells.Item($r, 3).Value = 69.900000000000006
}

# Remove the remaining sample rows (13-80): clear Serial/pdv, keep the
# empty, styled "valor" cell just like the untouched tail of the sheet.
for ($r = 13; $r -le 80; $r++) {
    $ws.Cells.Item($r, 1).Clear()
    $ws.Cells.Item($r, 2).Clear()
    $ws.Cells.Item($r, 3).ClearContents()
}

# Restore the selection to match the saved view state
[void]$ws.Range("A12:C12").Select()
